$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A duplicate/stale imported shift record was removed: row 466
# (Date 2025-04-10, Employee Daisy, Shift PH). Deleting the row
# shifts all subsequent rows up by one, which also makes the
# freshly-imported shift time override the previously customized one.
$ws.Rows.Item(466).Delete()
